# Updates cryptos list prices / volume(1h) figures, fixes a couple of
# mis-ordered rows (19/20 and 49/50/51), in line with the
# "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.658.89"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "1.582.81"
$ws.Range("E3").Value = "  -3.19%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.503"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.20%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.33"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.254"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0591"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0869"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("D12").Value = "1.808.27"
$ws.Range("E12").Value = "  -3.10%  "
$ws.Range("D13").Value = "1.589.49"
$ws.Range("E13").Value = "  -2.96%  "
$ws.Range("E14").Value = "  -4.05%  "
$ws.Range("E15").Value = "  -5.29%  "
$ws.Range("D16").Value = "27.637.48"
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.97%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0$([char]0x2083)0693"
$ws.Range("E19").Value = "  -3.59%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.62%  "
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("E27").Value = "  -2.51%  "
$ws.Range("E29").Value = "  -4.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.01%  "
$ws.Range("E31").Value = "  -3.54%  "
$ws.Range("E32").Value = "  -5.28%  "
$ws.Range("D33").Value = "1.375.50"
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("E34").Value = "  -5.23%  "
$ws.Range("E35").Value = "  -5.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.967"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.97%  "
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0166"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.540"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.70%  "
$ws.Range("E40").Value = "  -3.84%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  -3.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.26%  "
$ws.Range("D47").Value = "1.719.75"
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0$([char]0x2087)0999"
$ws.Range("E49").Value = "  -2.90%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0975"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.66%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0498"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.59%  "
